# Fixed LHS sampling to only sample across uncertainties (X) that vary
# (Ls still vary for all Ls) and rebuilt templates with
# PFLO:ALL_NO_STOPPING_DEFORESTATION_PLUR
#
# Concretely:
#   - Rename worksheet "strategy_id-5008" -> "strategy_id-5007"
#   - Add a new worksheet "strategy_id-5009" at the end, whose contents are
#     an exact duplicate of the (renamed) "strategy_id-5007" sheet.

$wb = $excel.ActiveWorkbook

# Rename strategy_id-5008 -> strategy_id-5007
$ws507 = $wb.Worksheets.Item("strategy_id-5008")
$ws507.Name = "strategy_id-5007"

# Duplicate it, placing the copy right after itself (i.e. at the end of
# the workbook), then rename the new copy to strategy_id-5009.
$ws507.Copy($null, $ws507)
$ws509 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws509.Name = "strategy_id-5009"
